$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# A1 sheet: replace the single evidence row, drop the old third row, add a
# basic page setup, and move the saved selection to A10.
# ---------------------------------------------------------------------------
$wsA1 = $wb.Worksheets.Item("A1")
$wsA1.Range("A2").Value = "C372AD12E2C95D0CF271A37A2B6F73D27737FBAC124F37C50373489F94D93D56"
$wsA1.Range("B2").Value = "birkoffGoNnft"
$wsA1.Rows.Item(3).Delete()
$wsA1.PageSetup.PaperSize = 9
$wsA1.PageSetup.Orientation = 1
$wsA1.Range("A10").Select()

# ---------------------------------------------------------------------------
# A2 sheet: replace the two evidence rows and append a third one, then select
# the three full data rows.
# ---------------------------------------------------------------------------
$wsA2 = $wb.Worksheets.Item("A2")
$wsA2.Range("A2").Value = "7A344EDB75ECFCE187A71FCC8F9B4DB43F0C149D654D3D6821CC41DE80CBBEC1"
$wsA2.Range("B2").Value = "birkoffGoNnft"
$wsA2.Range("C2").Value = "birkoff001"
$wsA2.Range("A3").Value = "47814F4B1B4BD28230997BEED0F7FB1FA668D94C12F2BBEFC015512C9AC703F7"
$wsA2.Range("B3").Value = "birkoffGoNnft"
$wsA2.Range("C3").Value = "birkoff002"
$wsA2.Range("B3:C3").Copy()
$wsA2.Range("B4:C4").PasteSpecial(-4122)
$wsA2.Range("A3:C3").Copy()
$wsA2.Range("A4:C4").PasteSpecial(-4122)
$wsA2.Range("A4").Value = "464D91C4A2FAA94244C29D60F94E5E106288CA0161DC8605B698ED8A63C250EB"
$wsA2.Range("B4").Value = "birkoffGoNnft"
$wsA2.Range("C4").Value = "birkoff003"
$wsA2.Range("A2:A4").EntireRow.Select()

# ---------------------------------------------------------------------------
# A3 sheet: fill in the real evidence values for the single data row.
# ---------------------------------------------------------------------------
$wsA3 = $wb.Worksheets.Item("A3")
$wsA3.Range("A2").Value = "68E5FEB2B0BD4630131E7C0808E9D7A45EC09F3563D7D1A3DF9A8626DC55EB36"
$wsA3.Range("B2").Value = "wasm.juno19us6395gfz2ehej6yj4hzv52zzpmwt55xy6e6xapd2a8lp3twltqpprnz7"
$wsA3.Range("C2").Value = "birkoff001"
$wsA3.Range("D2").Value = "uni-6"
$wsA3.Range("A2").EntireRow.Select()

# ---------------------------------------------------------------------------
# A4 sheet: fill in the real evidence values for the single data row.
# ---------------------------------------------------------------------------
$wsA4 = $wb.Worksheets.Item("A4")
$wsA4.Range("A2").Value = "0AB1D6CEE20724E8A552371CA2F4576F232C12D16E4AEDEBA83331126A8CC4F8"
$wsA4.Range("B2").Value = "ibc/145C27B96C1C9E3111F1B3602A56D8BD52BC6808E5A5F5BF60627C1D1D9E72B5"
$wsA4.Range("C2").Value = "birkoff002"
$wsA4.Range("D2").Value = "uptick_7000-2"
$wsA4.Range("B4").Select()

# ---------------------------------------------------------------------------
# A5 sheet: the evidence row is retyped from scratch (clears the custom row
# height), reusing A3's row format, then filled with the real values.
# ---------------------------------------------------------------------------
$wsA3.Range("A2:D2").Copy()
$wsA5 = $wb.Worksheets.Item("A5")
$wsA5.Rows.Item(2).Delete()
$wsA5.Range("A2:D2").PasteSpecial(-4122)
$wsA5.Range("A2").Value = "E6567D765DADE5C28C3253F82051E896CB13658797305BCB3F5C24679E74B85D"
$wsA5.Range("B2").Value = "wasm.juno19us6395gfz2ehej6yj4hzv52zzpmwt55xy6e6xapd2a8lp3twltqpprnz7"
$wsA5.Range("C2").Value = "birkoff001"
$wsA5.Range("D2").Value = "uni-6"
$wsA5.Range("B7").Select()

# ---------------------------------------------------------------------------
# A6 sheet: same treatment as A5, reusing A4's row format.
# ---------------------------------------------------------------------------
$wsA4.Range("A2:D2").Copy()
$wsA6 = $wb.Worksheets.Item("A6")
$wsA6.Rows.Item(2).Delete()
$wsA6.Range("A2:D2").PasteSpecial(-4122)
$wsA6.Range("A2").Value = "B3DAD9026D48303ACA9F519D9F9EEF5AA927D6CDAC63E00983709E9245442428"
$wsA6.Range("B2").Value = "ibc/145C27B96C1C9E3111F1B3602A56D8BD52BC6808E5A5F5BF60627C1D1D9E72B5"
$wsA6.Range("C2").Value = "birkoff002"
$wsA6.Range("D2").Value = "uptick_7000-2"
$wsA6.Range("B12").Select()

# ---------------------------------------------------------------------------
# Finally, move the active tab/selection to the Info sheet (this also clears
# tabSelected on whichever sheet previously held it).
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Activate()
$wsInfo.Range("E2").Select()
